$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append the new mail-log row (row 29) ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A29").Value = "Vragen over handleiding"
$logs.Range("B29").Value = "mailmind.test@zohomail.eu"
$logs.Range("C29").Value = "Waar kan ik de handleiding van product X vinden?"
$logs.Range("D29").Value = "Productinformatie"
$logs.Range("E29").Value = "Beste klant,`nBedankt voor uw vraag. De handleiding van product X is te vinden op onze website onder de rubriek 'Ondersteuning' of 'Downloads'. Mocht u hier hulp bij nodig hebben, aarzel dan niet om contact met ons op te nemen.`nMet vriendelijke groet,`n[Naam] - Klantenservice"
$logs.Range("F29").Value = "2025-06-22 18:59:16"
$logs.Range("G29").Value = "Ja"

# Excel auto-sizes the row for the wrapped multi-line text in E29; put it
# back to the sheet's normal (non-custom) row height like the other rows.
$logs.Rows.Item(29).AutoFit() | Out-Null

# --- Sheet "Dashboard": refresh the category summary table (rows 3-5) ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Productinformatie"
$dash.Range("B3").Value = 4

$dash.Range("A4").Value = "Afmelding / Nieuwsbrief"
$dash.Range("B4").Value = 3

$dash.Range("A5").Value = "Retour / Terugbetaling"
$dash.Range("B5").Value = 3

# --- extend the conditional-formatting ranges on "Logs" to cover the new row ---
$logs.Range("D2:D29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D29"))
$logs.Range("G2:G29").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G29"))
